$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.408.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.19%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.587.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.20%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +1.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'213.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.11%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.23%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +1.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'24.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.91%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.57%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.15%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.814.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.592.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.72%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.98%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.72%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'28.420.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'63.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'229.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.90%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.51%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +1.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.52%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.39%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.04%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.81%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +1.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.38%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.20%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.34%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.398.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.60%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -8.90%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0167"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.64%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D40").Value = "'0.542"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.15%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.42%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.06%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.99%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.981"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.82%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'64.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.73%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.726.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.43%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'87.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.61%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +4.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.31%  "
$ws.Range("E51").Style = "Normal"
